$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the row 14 totals: each column now sums only rows 12:13 instead of 12:16
$ws.Range("B14").Formula = "=SUM(B12:B13)"
$ws.Range("C14").Formula = "=SUM(C12:C13)"
$ws.Range("D14:N14").FormulaR1C1 = "=SUM(R[-2]C:R[-1]C)"

# Move the active selection to N17
$ws.Range("N17").Select()
